$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 2 for the new leaderboard entry
# "UI-TARS-1.5 (100 steps)" — this shifts every existing row down by one.
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the UI-TARS-1.5 entry's data.
$ws.Range("A2").Value = "UI-TARS-1.5 (100 steps)"
$ws.Range("B2").Value = "ByteDance Seed & Tsinghua University"
$ws.Range("C2").Value = "https://seed-tars.com/1.5"
$ws.Range("D2").Value = "Qin et al., '24"
$ws.Range("E2").Value = "—"
$ws.Range("F2").Value = 42.5
$ws.Range("G2").Value = "Apr 17, 2025"

# Rebuild the hyperlinks collection: the old single hyperlink (Kimi-VL's
# paper link, previously anchored at C27) now lives at C28 after the row
# insert, and we add the new hyperlink for the UI-TARS-1.5 paper link at C2.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C28"), "https://arxiv.org/abs/2504.07491", [System.Type]::Missing, [System.Type]::Missing, "https://arxiv.org/abs/2504.07491")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://seed-tars.com/1.5", [System.Type]::Missing, [System.Type]::Missing, "https://seed-tars.com/1.5")

# Apply the standard "Hyperlink" cell style (matching the existing paper-link
# cells elsewhere in the column) to both linked cells.
$ws.Range("C2").Style = "Hyperlink"
$ws.Range("C28").Style = "Hyperlink"

# Move the selection/view to match the saved workbook state (top-left reset,
# active cell G5).
$ws.Activate()
$ws.Range("G5").Select()
